$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column E slightly (stored OOXML width 23 -> 24; ColumnWidth is in
# character units and the engine adds the standard ~0.8333 padding offset
# when serializing, so back the padding out to land exactly on 24)
$ws.Columns.Item(5).ColumnWidth = 23.166666666666668

# Row 3 (PORCELANATO): updated VENTA / POR CUMPLIR / CUMPLIMIENTO
$ws.Range("D3").Value = 13968.27
$ws.Range("E3").Value = -244.9300000000003
$ws.Range("F3").Value = 1.017847695969057

# Row 4 (TOTAL): updated VENTA / POR CUMPLIR / CUMPLIMIENTO
$ws.Range("D4").Value = 19564.77
$ws.Range("E4").Value = -5841.43
$ws.Range("F4").Value = 1.425656582143997
